$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("K1").Value = "testbench"
$ws.Range("K2").Value = "SYS-110.tbc"
$ws.Range("K3").Value = "SYS-110.tbc"
$ws.Range("K4").Value = "SYS-110.tbc"

$ws.Range("L3").Select()
